# Apply updated dSF (column F) values, re-pulled from source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new column F value
$updates = @{
    2  = -4
    3  = -4
    6  = -3
    7  = 3
    8  = 2
    12 = 1
    13 = -1
    14 = 3
    16 = 2
    29 = -1
    30 = -2
    31 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
